$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.194.26"
$ws.Range("E2").Value = "  +2.98%  "
$ws.Range("D3").Value = "1.822.56"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("D4").Value = "0.9967"
$ws.Range("E4").Value = "  -0.72%  "
$ws.Range("D5").Value = "338.92"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("D6").Value = "0.9930"
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("D7").Value = "0.3933"
$ws.Range("E7").Value = "  +3.25%  "
$ws.Range("D8").Value = "0.3503"
$ws.Range("E8").Value = "  +0.89%  "
$ws.Range("D9").Value = "48.43"
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").Value = "1.205"
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("D11").Value = "0.07611"
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("D12").Value = "0.9945"
$ws.Range("E12").Value = "  -0.74%  "
$ws.Range("D13").Value = "22.27"
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("D14").Value = "6.563"
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("D15").Value = "1.817.77"
$ws.Range("E15").Value = "  +1.22%  "
$ws.Range("D16").Value = "7.224"
$ws.Range("E16").Value = "  +2.26%  "
$ws.Range("D17").Value = "0.00001110"
$ws.Range("E17").Value = "  +0.91%  "
$ws.Range("D18").Value = "0.06711"
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("D19").Value = "85.53"
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("D20").Value = "0.9945"
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("D21").Value = "17.93"
$ws.Range("E21").Value = "  +3.12%  "
$ws.Range("D22").Value = "6.593"
$ws.Range("E22").Value = "  +1.01%  "
$ws.Range("D23").Value = "28.191.82"
$ws.Range("E23").Value = "  +2.97%  "
$ws.Range("D24").Value = "12.86"
$ws.Range("E24").Value = "  +2.39%  "
$ws.Range("D25").Value = "2.398"
$ws.Range("E25").Value = "  -1.65%  "
$ws.Range("D26").Value = "1.547"
$ws.Range("E26").Value = "  +3.17%  "
$ws.Range("D27").Value = "2.589"
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("D28").Value = "21.52"
$ws.Range("E28").Value = "  +0.41%  "
$ws.Range("D29").Value = "154.92"
$ws.Range("E29").Value = "  +1.51%  "
$ws.Range("D30").Value = "2.024.77"
$ws.Range("E30").Value = "  +1.32%  "
$ws.Range("D31").Value = "136.15"
$ws.Range("E31").Value = "  +1.25%  "
$ws.Range("D32").Value = "6.216"
$ws.Range("E32").Value = "  +1.09%  "
$ws.Range("D33").Value = "4.034"
$ws.Range("E33").Value = "  -0.63%  "
$ws.Range("D34").Value = "0.08853"
$ws.Range("E34").Value = "  +1.62%  "
$ws.Range("D35").Value = "13.40"
$ws.Range("E35").Value = "  +0.93%  "
$ws.Range("D36").Value = "5.567"
$ws.Range("E36").Value = "  +2.12%  "
$ws.Range("D39").Value = "0.06574"
$ws.Range("E39").Value = "  +2.93%  "
$ws.Range("D40").Value = "1.615"
$ws.Range("E40").Value = "  -5.04%  "
$ws.Range("D41").Value = "0.2239"
$ws.Range("E41").Value = "  +1.12%  "
$ws.Range("D42").Value = "1.271"
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("D43").Value = "8.584"
$ws.Range("E43").Value = "  -4.26%  "
$ws.Range("D44").Value = "14.69"
$ws.Range("E44").Value = "  +1.21%  "
$ws.Range("D45").Value = "0.6538"
$ws.Range("E45").Value = "  +0.98%  "
$ws.Range("D46").Value = "3.879"
$ws.Range("E46").Value = "  +0.24%  "
$ws.Range("D47").Value = "2.178"
$ws.Range("E47").Value = "  +2.33%  "
$ws.Range("D48").Value = "132.76"
$ws.Range("E48").Value = "  +1.75%  "
$ws.Range("D49").Value = "0.07231"
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").Value = "80.82"
$ws.Range("E50").Value = "  +1.63%  "
$ws.Range("D51").Value = "1.163"
$ws.Range("E51").Value = "  +3.41%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.02452"
$ws.Range("E37").Value = "  +4.82%  "

$ws.Range("B38").Value = "TheSandbox"
$ws.Range("C38").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D38").Value = "0.6983"
$ws.Range("E38").Value = "  +1.04%  "
